$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 295 (shifts the existing rows 295-362 down to 296-363,
# which is exactly what the target diff shows: every old row N (295<=N<=362)
# reappears unchanged as row N+1, and the sheet's dimension grows from
# A1:R362 to A1:R363).
$ws.Rows("295").Insert()

# Populate the newly inserted row 295 with the new record.
$ws.Range("A295").Value = 10
$ws.Range("B295").Value = "Vega Modelo de Temuco"
$ws.Range("C295").Value = "La Araucanía"
$ws.Range("D295").Value = 45211
$ws.Range("E295").Value = 9
$ws.Range("F295").Value = 100112013
$ws.Range("G295").Value = "Alcachofa"
$ws.Range("H295").Value = "Española"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 240
$ws.Range("K295").Value = 12000
$ws.Range("L295").Value = 13000
$ws.Range("M295").Value = 12583
$ws.Range("N295").Value = "`$/caja 30 unidades"
$ws.Range("O295").Value = "Región Metropolitana"
$ws.Range("P295").Value = 419
$ws.Range("Q295").Value = 30
$ws.Range("R295").Value = "Hortaliza"
